$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row tweak: I1 "Python Predict Time " -> "Rust Predict Time " ---
$ws.Range("I1").Value = "Rust Predict Time "

# --- Row 2 (California Housing Dataset): H2 "9.515 ms" -> "0.497 ms" ---
$ws.Range("H2").Value = "0.497 ms"

# --- Row 3: Boston Housing Dataset ---
$ws.Range("A3").Value = "Boston Housing Dataset"
$ws.Range("B3").Value = -0.1779
$ws.Range("C3").Value = 0.58426999999999996
$ws.Range("D3").Value = 31.671299999999999
$ws.Range("E3").Value = 4466072770.7544003
$ws.Range("F3").Value = "1.733 ms"
$ws.Range("G3").Value = "24.667 ms"
$ws.Range("H3").Value = " 0.343 ms"
$ws.Range("I3").Value = "1.380 ms"
$ws.Range("J3").Value = "Kaggle"
$ws.Range("K3").Value = "same dataset ,same spilt "

# --- Row 4: Auto MPG Dataset ---
$ws.Range("A4").Value = "Auto MPG Dataset"
$ws.Range("B4").Value = 0.212492825720139
$ws.Range("C4").Value = 0.57030999999999998
$ws.Range("D4").Value = 28.2358224073324
$ws.Range("E4").Value = 14.549307000000001
$ws.Range("F4").Value = "2.064 ms"
$ws.Range("G4").Value = "0.841 ms"
$ws.Range("H4").Value = "0.274 ms"
$ws.Range("I4").Value = "0.025 ms"

# --- Row 5: Advertising Dataset ---
$ws.Range("A5").Value = "Advertising Dataset"
$ws.Range("B5").Value = 0.90247132416945797
$ws.Range("C5").Value = 0.87556299999999998
$ws.Range("D5").Value = 2.6369407160816101
$ws.Range("E5").Value = 3.0560689999999999
$ws.Range("F5").Value = "1.298 ms"
$ws.Range("G5").Value = "0.331 ms"
$ws.Range("H5").Value = "0.245 ms"
$ws.Range("I5").Value = "0.028 ms"

# --- Row 6: Diabetes Dataset ---
$ws.Range("A6").Value = "Diabetes Dataset"
$ws.Range("B6").Value = 0.299000777026311
$ws.Range("C6").Value = 0.25683299999999998
$ws.Range("D6").Value = 0.16094369915212201
$ws.Range("E6").Value = 0.169215
$ws.Range("F6").Value = "1.548 ms"
$ws.Range("G6").Value = "1.191 ms"
$ws.Range("H6").Value = "0.252 ms"
$ws.Range("I6").Value = "0.054 ms"

# --- Row 7: Stock Market Dataset ---
$ws.Range("A7").Value = "Stock Market Dataset"
$ws.Range("B7").Value = 0.32612154927842202
$ws.Range("C7").Value = -0.20691899999999999
$ws.Range("D7").Value = 5070387311095220
$ws.Range("E7").Value = 7412450837751620
$ws.Range("F7").Value = "2.584 ms"
$ws.Range("G7").Value = "1.267 ms"
$ws.Range("H7").Value = "0.512 ms"
$ws.Range("I7").Value = "0.098 ms"

# --- Match the final selection recorded in the workbook (G7) ---
[void]$ws.Range("G7").Select()
